$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New value next to row 2 (F2 = "admin")
$ws.Range("F2").Value = "admin"

# New permission row (row 10)
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "HasPermissionToPublishEvent"

# Match the vertical-center style used by the rest of the A/B column cells
$ws.Range("A10:B10").VerticalAlignment = -4108

# Column B grows to fit the new (longer) permission name
$ws.Columns("B").ColumnWidth = 24.6

# Move the selection to the newly added cell
[void]$ws.Range("B10").Select()
